$d = $word.ActiveDocument

# Locate the paragraph that holds the Jinja tag:
#   {% if not staffReason == "N/A" %}
# (curly/smart quotes around N/A, exactly like the template uses elsewhere)
$quote1 = [char]0x201C
$quote2 = [char]0x201D
$needleWhole = "not staffReason == " + $quote1 + "N/A" + $quote2 + " "

$targetIndex = -1
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text -like "*$needleWhole*") {
        $targetIndex = $i
    }
}

if ($targetIndex -ne -1) {
    $p = $d.Paragraphs.Item($targetIndex)
    $pStart = $p.Range.Start
    $pText = $p.Range.Text

    $notStart = $pText.IndexOf("not ")
    $notEnd = $notStart + 4

    $condStart = $pText.IndexOf(" == " + $quote1)
    $condEnd = $pText.IndexOf("%}")

    # Remove the trailing ' == "N/A" ' portion (before the closing '%}') first,
    # so the earlier offsets for 'not ' stay valid.
    $r1 = $d.Range($pStart + $condStart, $pStart + $condEnd)
    $r1.Delete()

    # Remove the leading 'not ' portion.
    $r2 = $d.Range($pStart + $notStart, $pStart + $notEnd)
    $r2.Delete()

    # Rename staffReason -> staffHasReason in place (keeps the surrounding
    # proofErr spell-check markers and run formatting intact).
    $p2 = $d.Paragraphs.Item($targetIndex)
    $pStart2 = $p2.Range.Start
    $p2Text = $p2.Range.Text
    $nameStart = $p2Text.IndexOf("staffReason")
    $nameEnd = $nameStart + "staffReason".Length
    $r3 = $d.Range($pStart2 + $nameStart, $pStart2 + $nameEnd)
    $r3.Text = "staffHasReason"

    # Shrink the whole line (all runs) from 11pt (sz/szCs 22) down to 7pt
    # (sz/szCs 14) to visually hide/de-emphasize this Jinja control line.
    $p3 = $d.Paragraphs.Item($targetIndex)
    $p3.Range.Font.Size = 7
    $p3.Range.Font.SizeBi = 7
}
